$d = $word.ActiveDocument
try {
  $n = $d.CustomXMLParts.Count
  Write-Output "Count=$n"
} catch {
  Write-Output "ERR: $_"
}
try {
  $parts = $d.CustomXMLParts
  Write-Output "parts type got"
  $sel = $parts.SelectByNamespace("urn:microsoft-dynamics-nav/reports/BBC_WOSF_Sales_Credit_Memo/50201/")
  Write-Output "sel count = $($sel.Count)"
} catch {
  Write-Output "ERR2: $_"
}
